$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Add new headers for the Haunted Carriage event timers (columns K, L, M)
$ws.Range("K1").Value = "Haunted_Carriage_12PM"
$ws.Range("L1").Value = "Haunted_Carriage_830PM"
$ws.Range("M1").Value = "Haunted_Carriage_10PM"

# Size the new columns to fit their header text (matches the other bestFit columns)
$ws.Columns.Item(11).ColumnWidth = 22.25
$ws.Columns.Item(12).ColumnWidth = 23.25
$ws.Columns.Item(13).ColumnWidth = 22.25

# Reflect the new selection state (E8:E9) recorded in the saved view
$ws.Range("E8:E9").Select()
